# Life Point Prefab Pickup Added
# Update the ENTITY_PATH column (B) for the existing loot rows to point at
# the new shared prefab path, and append a new row describing the
# LIFE_POINT_1 pickup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loot_All")

# New Life Point pickup row's entity path (added first so it lands right
# after the existing shared strings in the table).
$ws.Cells.Item(9, 2).Value2 = "Assets/Prefabs/Pickups/LifePoint.prefab"

# Every existing data row (2-8) referenced the old
# "Assets/Scripts/LootTable/Loot/L_Heart_1.asset" entity path; it is now
# replaced by the shared key prefab path.
$newEntityPath = "Assets/Prefabs/Pickups/Key.prefab"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $newEntityPath
}

# Finish filling in the new Life Point pickup row.
$ws.Cells.Item(9, 1).Value2 = "LIFE_POINT_1"
$ws.Cells.Item(9, 3).Value2 = "Minor Heal"
$ws.Cells.Item(9, 4).Value2 = "COMMON"
